# Refresh the "cryptos" market-data table (price + 1h volume change).
# Two rows (50/51) also swap which coin (Aave/Cronos) occupies them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.696.88"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.699.67"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3931"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4037"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08849"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.398"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.133"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001321"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "1.707.20"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07044"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.091"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.78%  "
$ws.Range("D24").Value = "24.691.18"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.369"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.748"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +17.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.184"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09029"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.062"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.979"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2753"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02766"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09127"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.463"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7666"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7166"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.558"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.208"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.332"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07983"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07983"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.84%  "
